$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 838.2857
$ws.Range("I41").Value = 859.55554
$ws.Range("J41").Value = 800
$ws.Range("K41").Value = 859.55554
$ws.Range("L41").Value = 800
$ws.Range("M41").Value = -419.55554
$ws.Range("N41").Value = -1680

$ws.Range("H99").Value = 1025.9166
$ws.Range("I99").Value = 1025.9166
$ws.Range("K99").Value = 3077.7498
$ws.Range("M99").Value = -1579.7498

$ws.Range("H121").Value = 855.8125
$ws.Range("I121").Value = 548.4286
$ws.Range("J121").Value = 1094.8889
$ws.Range("K121").Value = 1645.2858
$ws.Range("L121").Value = 3284.6667
$ws.Range("M121").Value = 101.7142000000001
$ws.Range("N121").Value = -6778.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 385296.25
$ws.Range("I32").Value = 427490.16
$ws.Range("K32").Value = 427490.16
$ws.Range("M32").Value = -427203.16

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H61").Value = 346743.34
$ws.Range("I61").Value = 1728.625
$ws.Range("J61").Value = 2002814
$ws.Range("K61").Value = 1728.625
$ws.Range("L61").Value = 2002814
$ws.Range("M61").Value = -1516.625
$ws.Range("N61").Value = -2003238

$ws.Range("H74").Value = 3721.1892
$ws.Range("I74").Value = 907.0454999999999
$ws.Range("J74").Value = 7848.6
$ws.Range("K74").Value = 907.0454999999999
$ws.Range("L74").Value = 7848.6
$ws.Range("M74").Value = -33.04549999999995
$ws.Range("N74").Value = -9596.6

$ws.Range("H77").Value = 3721.1892
$ws.Range("I77").Value = 907.0454999999999
$ws.Range("J77").Value = 7848.6
$ws.Range("K77").Value = 4535.2275
$ws.Range("L77").Value = 39243
$ws.Range("M77").Value = -167.2275
$ws.Range("N77").Value = -47979

$ws.Range("H136").Value = 346743.34
$ws.Range("I136").Value = 1728.625
$ws.Range("J136").Value = 2002814
$ws.Range("K136").Value = 5185.875
$ws.Range("L136").Value = 6008442
$ws.Range("M136").Value = -2635.875
$ws.Range("N136").Value = -6013542

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3040.9666
$ws.Range("I86").Value = 2562.3635
$ws.Range("J86").Value = 3318.0527
$ws.Range("K86").Value = 2562.3635
$ws.Range("L86").Value = 3318.0527
$ws.Range("M86").Value = -1439.3635
$ws.Range("N86").Value = -5564.0527

$ws.Range("H89").Value = 3040.9666
$ws.Range("I89").Value = 2562.3635
$ws.Range("J89").Value = 3318.0527
$ws.Range("K89").Value = 12811.8175
$ws.Range("L89").Value = 16590.2635
$ws.Range("M89").Value = -7195.817499999999
$ws.Range("N89").Value = -27822.2635

$ws.Range("H134").Value = 1547.9231
$ws.Range("I134").Value = 1200.2222
$ws.Range("J134").Value = 2330.25
$ws.Range("K134").Value = 3600.6666
$ws.Range("L134").Value = 6990.75
$ws.Range("M134").Value = -1065.6666
$ws.Range("N134").Value = -12060.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30000.428
$ws.Range("J64").Value = 30000.428
$ws.Range("L64").Value = 30000.428
$ws.Range("N64").Value = -30496.428

$ws.Range("H67").Value = 30000.428
$ws.Range("J67").Value = 30000.428
$ws.Range("L67").Value = 30000.428
$ws.Range("N67").Value = -31716.428

$ws.Range("H95").Value = 13958.875
$ws.Range("J95").Value = 13958.875
$ws.Range("L95").Value = 13958.875
$ws.Range("N95").Value = -19450.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1498.75
$ws.Range("I68").Value = 725.8182
$ws.Range("J68").Value = 10001
$ws.Range("K68").Value = 2177.4546
$ws.Range("L68").Value = 30003
$ws.Range("M68").Value = -1366.4546
$ws.Range("N68").Value = -31625

$ws.Range("H71").Value = 1498.75
$ws.Range("I71").Value = 725.8182
$ws.Range("J71").Value = 10001
$ws.Range("K71").Value = 6532.3638
$ws.Range("L71").Value = 90009
$ws.Range("M71").Value = -2476.3638
$ws.Range("N71").Value = -98121

$ws.Range("H121").Value = 5112.3335
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 5613.674
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 16841.022
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -19461.022

$ws.Range("H132").Value = 56606.277
$ws.Range("I132").Value = 913.5
$ws.Range("J132").Value = 101160.5
$ws.Range("K132").Value = 8221.5
$ws.Range("L132").Value = 910444.5
$ws.Range("M132").Value = -5691.5
$ws.Range("N132").Value = -915504.5

$ws.Range("H136").Value = 3704.125
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 3947.5715
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 11842.7145
$ws.Range("M136").Value = -900
$ws.Range("N136").Value = -22042.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2892.4614
$ws.Range("I80").Value = 2266.2856
$ws.Range("J80").Value = 3623
$ws.Range("K80").Value = 2266.2856
$ws.Range("L80").Value = 3623
$ws.Range("M80").Value = -1268.2856
$ws.Range("N80").Value = -5619

$ws.Range("H83").Value = 2892.4614
$ws.Range("I83").Value = 2266.2856
$ws.Range("J83").Value = 3623
$ws.Range("K83").Value = 11331.428
$ws.Range("L83").Value = 18115
$ws.Range("M83").Value = -6339.428
$ws.Range("N83").Value = -28099

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 39333.332
$ws.Range("J64").Value = 39333.332
$ws.Range("L64").Value = 39333.332
$ws.Range("N64").Value = -39783.332

$ws.Range("H67").Value = 39333.332
$ws.Range("J67").Value = 39333.332
$ws.Range("L67").Value = 39333.332
$ws.Range("N67").Value = -40893.332

$ws.Range("H68").Value = 1744.4103
$ws.Range("I68").Value = 1633.4166
$ws.Range("J68").Value = 1922
$ws.Range("K68").Value = 1633.4166
$ws.Range("L68").Value = 1922
$ws.Range("M68").Value = -884.4166
$ws.Range("N68").Value = -3420

$ws.Range("H71").Value = 1744.4103
$ws.Range("I71").Value = 1633.4166
$ws.Range("J71").Value = 1922
$ws.Range("K71").Value = 8167.083000000001
$ws.Range("L71").Value = 9610
$ws.Range("M71").Value = -4423.083000000001
$ws.Range("N71").Value = -17098

$ws.Range("H74").Value = 48000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 48000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 48000
$ws.Range("N74").Value = -49996
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 48000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 48000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 144000
$ws.Range("N77").Value = -153984
$ws.Range("M77").ClearContents()

$ws.Range("H82").Value = 1317.5714
$ws.Range("I82").Value = 1078
$ws.Range("J82").Value = 1535.3636
$ws.Range("K82").Value = 1078
$ws.Range("L82").Value = 1535.3636
$ws.Range("M82").Value = -717
$ws.Range("N82").Value = -2257.3636

$ws.Range("H85").Value = 1317.5714
$ws.Range("I85").Value = 1078
$ws.Range("J85").Value = 1535.3636
$ws.Range("K85").Value = 1078
$ws.Range("L85").Value = 1535.3636
$ws.Range("M85").Value = 170
$ws.Range("N85").Value = -4031.3636

$ws.Range("H122").Value = 41385.08
$ws.Range("I122").Value = 68917.13
$ws.Range("J122").Value = 3841.3635
$ws.Range("K122").Value = 206751.39
$ws.Range("L122").Value = 11524.0905
$ws.Range("M122").Value = -204301.39
$ws.Range("N122").Value = -16424.0905

$ws.Range("H136").Value = 5377.3706
$ws.Range("I136").Value = 1326.4445
$ws.Range("J136").Value = 13479.223
$ws.Range("K136").Value = 3979.3335
$ws.Range("L136").Value = 40437.669
$ws.Range("M136").Value = -1429.3335
$ws.Range("N136").Value = -45537.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3051.111
$ws.Range("I62").Value = 2750
$ws.Range("K62").Value = 2750
$ws.Range("M62").Value = -2126

$ws.Range("H63").Value = 14300
$ws.Range("I63").Value = 5933.3335
$ws.Range("J63").Value = 22666.666
$ws.Range("K63").Value = 5933.3335
$ws.Range("L63").Value = 22666.666
$ws.Range("M63").Value = -5309.3335
$ws.Range("N63").Value = -23914.666

$ws.Range("H65").Value = 3051.111
$ws.Range("I65").Value = 2750
$ws.Range("K65").Value = 13750
$ws.Range("M65").Value = -10630

$ws.Range("H66").Value = 14300
$ws.Range("I66").Value = 5933.3335
$ws.Range("J66").Value = 22666.666
$ws.Range("K66").Value = 17800.0005
$ws.Range("L66").Value = 67999.99800000001
$ws.Range("M66").Value = -14680.0005
$ws.Range("N66").Value = -74239.99800000001
